$d = $word.ActiveDocument

$replacements = @(
    @{Old = "306×2=612"; New = "595×6=3570"},
    @{Old = "285×2=570"; New = "812×9=7308"},
    @{Old = "808×7=5656"; New = "656×6=3936"},
    @{Old = "581×8=4648"; New = "960×9=8640"},
    @{Old = "711×7=4977"; New = "826×4=3304"},
    @{Old = "728×5=3640"; New = "556×8=4448"},
    @{Old = "225×9=2025"; New = "753×7=5271"},
    @{Old = "743×4=2972"; New = "539×5=2695"},
    @{Old = "952×5=4760"; New = "730×5=3650"},
    @{Old = "196×2=392"; New = "178×7=1246"},
    @{Old = "115×5=575"; New = "338×9=3042"},
    @{Old = "564×7=3948"; New = "172×4=688"},
    @{Old = "290×6=1740"; New = "954×6=5724"},
    @{Old = "573×8=4584"; New = "846×9=7614"},
    @{Old = "688×3=2064"; New = "634×3=1902"},
    @{Old = "261×9=2349"; New = "299×9=2691"},
    @{Old = "213×7=1491"; New = "634×3=1902"},
    @{Old = "171×7=1197"; New = "977×9=8793"},
    @{Old = "578×8=4624"; New = "602×3=1806"},
    @{Old = "147×7=1029"; New = "498×3=1494"},
    @{Old = "847×8=6776"; New = "937×3=2811"},
    @{Old = "706×4=2824"; New = "663×3=1989"},
    @{Old = "268×2=536"; New = "243×2=486"},
    @{Old = "521×6=3126"; New = "361×5=1805"},
    @{Old = "394×4=1576"; New = "794×2=1588"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
